$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value2 = 200.1875
$ws.Range("I12").Value2 = 198.3
$ws.Range("J12").Value2 = 203.33333
$ws.Range("K12").Value2 = 198.3
$ws.Range("L12").Value2 = 203.33333
$ws.Range("M12").Value2 = -28.30000000000001
$ws.Range("N12").Value2 = -543.3333299999999
$ws.Range("H64").Value2 = 3895.72
$ws.Range("I64").Value2 = 3806.0833
$ws.Range("K64").Value2 = 3806.0833
$ws.Range("M64").Value2 = -3558.0833
$ws.Range("H67").Value2 = 3895.72
$ws.Range("I67").Value2 = 3806.0833
$ws.Range("K67").Value2 = 3806.0833
$ws.Range("M67").Value2 = -2948.0833
$ws.Range("H74").Value2 = 3500
$ws.Range("I74").Value2 = 3000
$ws.Range("K74").Value2 = 3000
$ws.Range("M74").Value2 = -2064
$ws.Range("H76").Value2 = 3485.2104
$ws.Range("I76").Value2 = 2997.6743
$ws.Range("J76").Value2 = 4120.485
$ws.Range("K76").Value2 = 2997.6743
$ws.Range("L76").Value2 = 4120.485
$ws.Range("M76").Value2 = -2682.6743
$ws.Range("N76").Value2 = -4750.485
$ws.Range("H77").Value2 = 3500
$ws.Range("I77").Value2 = 3000
$ws.Range("K77").Value2 = 15000
$ws.Range("M77").Value2 = -10320
$ws.Range("H79").Value2 = 3485.2104
$ws.Range("I79").Value2 = 2997.6743
$ws.Range("J79").Value2 = 4120.485
$ws.Range("K79").Value2 = 2997.6743
$ws.Range("L79").Value2 = 4120.485
$ws.Range("M79").Value2 = -1905.6743
$ws.Range("N79").Value2 = -6304.485
$ws.Range("H100").Value2 = 1623.7142
$ws.Range("I100").Value2 = 1750
$ws.Range("J100").Value2 = 1573.2
$ws.Range("K100").Value2 = 1750
$ws.Range("L100").Value2 = 1573.2
$ws.Range("M100").Value2 = -1209
$ws.Range("N100").Value2 = -2655.2
$ws.Range("H107").Value2 = 509.94446
$ws.Range("I107").Value2 = 479.9375
$ws.Range("J107").Value2 = 750
$ws.Range("K107").Value2 = 479.9375
$ws.Range("L107").Value2 = 750
$ws.Range("M107").Value2 = 1440.0625
$ws.Range("N107").Value2 = -4590
$ws.Range("H115").Value2 = 20000948
$ws.Range("I115").Value2 = 20000948
$ws.Range("K115").Value2 = 60002844
$ws.Range("M115").Value2 = -60001277
$ws.Range("H125").Value2 = 1193.5454
$ws.Range("I125").Value2 = 915.5
$ws.Range("K125").Value2 = 8239.5
$ws.Range("M125").Value2 = -5779.5
$ws.Range("H129").Value2 = 959.2967
$ws.Range("I129").Value2 = 572.46155
$ws.Range("J129").Value2 = 1023.7692
$ws.Range("K129").Value2 = 1717.38465
$ws.Range("L129").Value2 = 3071.3076
$ws.Range("M129").Value2 = 3282.61535
$ws.Range("N129").Value2 = -13071.3076
$ws.Range("H132").Value2 = 1846.1333
$ws.Range("I132").Value2 = 1892.5518
$ws.Range("J132").Value2 = 500
$ws.Range("K132").Value2 = 5677.6554
$ws.Range("L132").Value2 = 1500
$ws.Range("M132").Value2 = -3147.6554
$ws.Range("N132").Value2 = -6560
$ws.Range("H134").Value2 = 42000
$ws.Range("J134").Value2 = 42000
$ws.Range("L134").Value2 = 42000
$ws.Range("N134").Value2 = -52140
$ws.Range("H137").Value2 = 5705988
$ws.Range("I137").Value2 = 7143906.5
$ws.Range("K137").Value2 = 21431719.5
$ws.Range("M137").Value2 = -21429169.5
$ws.Range("H138").Value2 = 3154.549
$ws.Range("I138").Value2 = 689.0526
$ws.Range("J138").Value2 = 4618.4375
$ws.Range("K138").Value2 = 2067.1578
$ws.Range("L138").Value2 = 13855.3125
$ws.Range("M138").Value2 = 3072.8422
$ws.Range("N138").Value2 = -24135.3125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 4559.737
$ws.Range("I32").Value2 = 2667.3142
$ws.Range("J32").Value2 = 26638
$ws.Range("K32").Value2 = 2667.3142
$ws.Range("L32").Value2 = 26638
$ws.Range("M32").Value2 = -2380.3142
$ws.Range("N32").Value2 = -27212
$ws.Range("H110").Value2 = 4757.4707
$ws.Range("I110").Value2 = 3952.3635
$ws.Range("J110").Value2 = 6233.5
$ws.Range("K110").Value2 = 3952.3635
$ws.Range("L110").Value2 = 6233.5
$ws.Range("M110").Value2 = -1907.3635
$ws.Range("N110").Value2 = -10323.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value2 = 429.29413
$ws.Range("I64").Value2 = 507.125
$ws.Range("J64").Value2 = 360.1111
$ws.Range("K64").Value2 = 507.125
$ws.Range("L64").Value2 = 360.1111
$ws.Range("M64").Value2 = -282.125
$ws.Range("N64").Value2 = -810.1111000000001
$ws.Range("H67").Value2 = 429.29413
$ws.Range("I67").Value2 = 507.125
$ws.Range("J67").Value2 = 360.1111
$ws.Range("K67").Value2 = 507.125
$ws.Range("L67").Value2 = 360.1111
$ws.Range("M67").Value2 = 272.875
$ws.Range("N67").Value2 = -1920.1111
$ws.Range("H99").Value2 = 3207.8572
$ws.Range("I99").Value2 = 2702
$ws.Range("J99").Value2 = 3488.889
$ws.Range("K99").Value2 = 2702
$ws.Range("L99").Value2 = 3488.889
$ws.Range("M99").Value2 = -1204
$ws.Range("N99").Value2 = -6484.889

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value2 = 1105.2858
$ws.Range("I107").Value2 = 970.3333
$ws.Range("J107").Value2 = 1206.5
$ws.Range("K107").Value2 = 970.3333
$ws.Range("L107").Value2 = 1206.5
$ws.Range("M107").Value2 = 949.6667
$ws.Range("N107").Value2 = -5046.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value2 = 46994.91
$ws.Range("I4").Value2 = 436.83334
$ws.Range("J4").Value2 = 56971.645
$ws.Range("K4").Value2 = 1310.50002
$ws.Range("L4").Value2 = 170914.935
$ws.Range("M4").Value2 = -1198.50002
$ws.Range("N4").Value2 = -171138.935
$ws.Range("H92").Value2 = 555808.4399999999
$ws.Range("I92").Value2 = 1000185.8
$ws.Range("J92").Value2 = 336.75
$ws.Range("K92").Value2 = 3000557.4
$ws.Range("L92").Value2 = 1010.25
$ws.Range("M92").Value2 = -2999309.4
$ws.Range("N92").Value2 = -3506.25
$ws.Range("H97").Value2 = 943.3077
$ws.Range("I97").Value2 = 373.25
$ws.Range("J97").Value2 = 1196.6666
$ws.Range("K97").Value2 = 1119.75
$ws.Range("L97").Value2 = 3589.9998
$ws.Range("M97").Value2 = -623.75
$ws.Range("N97").Value2 = -4581.9998
$ws.Range("H113").Value2 = 1326753.1
$ws.Range("I113").Value2 = 5747545
$ws.Range("J113").Value2 = 515.6
$ws.Range("K113").Value2 = 17242635
$ws.Range("L113").Value2 = 1546.8
$ws.Range("M113").Value2 = -17240465
$ws.Range("N113").Value2 = -5886.8
$ws.Range("H122").Value2 = 790.9474
$ws.Range("I122").Value2 = 649.8889
$ws.Range("J122").Value2 = 917.9
$ws.Range("K122").Value2 = 5849.0001
$ws.Range("L122").Value2 = 8261.1
$ws.Range("M122").Value2 = -3399.0001
$ws.Range("N122").Value2 = -13161.1
$ws.Range("H131").Value2 = 799.2
$ws.Range("I131").Value2 = 382.22223
$ws.Range("K131").Value2 = 1146.66669
$ws.Range("M131").Value2 = 3893.33331

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value2 = 1474.9762
$ws.Range("I102").Value2 = 1452.5758
$ws.Range("K102").Value2 = 1452.5758
$ws.Range("M102").Value2 = 169.4241999999999
$ws.Range("H126").Value2 = 2190.2
$ws.Range("I126").Value2 = 1148.5714
$ws.Range("J126").Value2 = 3101.625
$ws.Range("K126").Value2 = 3445.7142
$ws.Range("L126").Value2 = 9304.875
$ws.Range("M126").Value2 = -975.7142000000003
$ws.Range("N126").Value2 = -14244.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value2 = 2218.9
$ws.Range("I7").Value2 = 2242
$ws.Range("J7").Value2 = 2165
$ws.Range("K7").Value2 = 2242
$ws.Range("L7").Value2 = 2165
$ws.Range("M7").Value2 = -2130
$ws.Range("N7").Value2 = -2389
$ws.Range("H61").Value2 = 3903.2646
$ws.Range("I61").Value2 = 4493.8276
$ws.Range("K61").Value2 = 4493.8276
$ws.Range("M61").Value2 = -4291.8276
$ws.Range("H113").Value2 = 3903.2646
$ws.Range("I113").Value2 = 4493.8276
$ws.Range("K113").Value2 = 4493.8276
$ws.Range("M113").Value2 = -2323.8276
$ws.Range("H126").Value2 = 2218.9
$ws.Range("I126").Value2 = 2242
$ws.Range("J126").Value2 = 2165
$ws.Range("K126").Value2 = 6726
$ws.Range("L126").Value2 = 6495
$ws.Range("M126").Value2 = -4256
$ws.Range("N126").Value2 = -11435
$ws.Range("H128").Value2 = 39564.5
$ws.Range("J128").Value2 = 39564.5
$ws.Range("L128").Value2 = 39564.5
$ws.Range("N128").Value2 = -49524.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value2 = 2527
$ws.Range("J96").Value2 = 2840.5
$ws.Range("L96").Value2 = 2840.5
$ws.Range("N96").Value2 = -5586.5
$ws.Range("H113").Value2 = 922.7778
$ws.Range("I113").Value2 = 773
$ws.Range("J113").Value2 = 1042.6
$ws.Range("K113").Value2 = 2319
$ws.Range("L113").Value2 = 3127.8
$ws.Range("M113").Value2 = -149
$ws.Range("N113").Value2 = -7467.799999999999
$ws.Range("H135").Value2 = 25358
$ws.Range("J135").Value2 = 25358
$ws.Range("L135").Value2 = 25358
$ws.Range("N135").Value2 = -35498
